$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the A-column "Marked" flag to 0 for the filtered-out rows
$rowsToClear = @(76, 77, 79, 80, 81, 85, 133)
foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 1).Value = 0
}

# Update the active selection to reflect the last reviewed cell
$ws.Range("A133").Select()
